$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 data for Donna Davis, reusing the shared Password/Access values from row 4
$ws.Range("A5").Value = "Donna Davis"
$ws.Range("B5").Value = "donna.davis@hgv.com"
$ws.Range("C5").Value = $ws.Range("C4").Value()
$ws.Range("D5").Value = $ws.Range("D4").Value()

# Match row height used by the other data rows
$ws.Rows.Item(5).RowHeight = 14.9

# Add the mailto hyperlink for the new email address
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:donna.davis@hgv.com", [System.Type]::Missing, [System.Type]::Missing, "donna.davis@hgv.com")

# Update the active selection to D5, matching the saved view state
$ws.Range("D5").Select()
